# Update column D ("Name of Algo" result column) values produced by the KNN
# imputation run for terrestrial_mammals / combination_3_ABCDF / D / 20 / seed4.
#
# Only five cells in column D change value; everything else in the sheet stays
# the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value  = -8.319000000000001
$ws.Range("D13").Value = -7.662999999999999
$ws.Range("D16").Value = -8.183
$ws.Range("D18").Value = -8.176
$ws.Range("D20").Value = -8.083000000000002
